# Apply the "Subida versión inicial de las consultas a la BBDD para las recetas" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) F55: update formula 1+1 -> 1+1+0.5 (value becomes 2.5)
$ws.Range("F55").Formula = "=1+1+0.5"

# 2) C55: fix stray white-fill style so it matches the other "pending" cells (C54/C58/C59)
$ws.Range("C55").Interior.Color = $ws.Range("C54").Interior.Color

# 3) F56: new formula cell 1+2 = 3
$ws.Range("F56").Formula = "=1+2"

# 4) Insert a new row at 57 (shifts old rows 57-64 down to 58-65)
$ws.Rows("57:57").Insert()

# Copy formatting from the row that is now at 58 (a "validation" style row,
# same visual style as the new row should have) into the new row 57
$ws.Range("B58").Copy($ws.Range("B57"))

# Now set the real content of the new row 57
$ws.Range("B57").Value2 = "Validar que la consulta sale por pantalla"
$ws.Range("D57").Value2 = "Sergio"
$ws.Range("E57").Value2 = 0.1

$wb.Save()
